$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts all existing columns
# (A:V) right by one (to B:W) and copies column A's formatting into B.
$ws.Columns("A:A").Insert()

# Header row (row 2): label the new column "Match ID"
$ws.Range("A2").Value = "Match ID"

# Data rows 4-19 (and the hidden total row 20) all belong to match id 4
for ($r = 4; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = 4
}

# Row 3 is a hidden blank separator row; leave A3 blank but give it the
# same (bold, borderless) style as the rest of column A.
$ws.Range("A2:A20").Font.Bold = $true
$ws.Range("A2:A20").Borders.LineStyle = -4142

# Row 20 (hidden total row) keeps the default (non-bold) style for A20
$ws.Cells.Item(20, 1).Font.Bold = $false

# Update the sheet selection to match the new layout
$ws.Range("A2:A19").Select()
